# Update "want to go" counts (column F) on the 展览 and 全部类型 sheets
# F2: 229 -> 230
# F3: 163 -> 164
# F5: 2   -> 3

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 230
    $ws.Range("F3").Value = 164
    $ws.Range("F5").Value = 3
}
